$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.108.12'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.17%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.830.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.27%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.38%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.12%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6265'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '

# Row 7
$ws.Range("E7").Value = '  +0.15%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07496'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.34%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2930'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07704'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.020'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.94%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.796.41'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.27%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6684'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.31%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.77'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.60%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009386'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.16%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.984'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.49%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.091.83'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.13%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.065.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.01%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.60'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.39%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '223.35'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.62%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '

# Row 24
$ws.Range("E24").Value = '  +0.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '160.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.23%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1400'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.98%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.504'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.27%  '

# Row 28
$ws.Range("E28").Value = '  -0.16%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.489'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.14%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05782'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.98%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.160'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.29%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.123'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.52%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.210'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.34%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7418'
$ws.Range("D34").Style = "Normal"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.830'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.83%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.139'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.05%  '

# Row 37
$ws.Range("E37").Value = '  -0.33%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.232.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.83%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.762'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01775'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.49%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.504'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.19%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8937'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.72%  '

# Row 43
$ws.Range("E43").Value = '  +0.08%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '102.26'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.86'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.50%  '

# Row 46
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.959.19'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.27%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5091'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.42%  '

# Row 48
$ws.Range("B48").Value = 'TheSandbox'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4068'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.54%  '

# Row 49
$ws.Range("B49").Value = 'XinFinNetwork'
$ws.Range("C49").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07484'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +12.12%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.017'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.51%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05831'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.16%  '
